{"js": "// Fix the typo \"I2C_HIGH_SPPED_MODE\" -> \"I2C_HIGH_SPEED_MODE\" in the\n// paragraph describing the I2C speed-mode change.\nconst body = context.document.body;\n\nconst results = body.search(\"I2C_HIGH_SPPED_MODE\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"I2C_HIGH_SPEED_MODE\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Fix the typo \"I2C_HIGH_SPPED_MODE\" -> \"I2C_HIGH_SPEED_MODE\" in the\n# paragraph describing the I2C speed-mode change.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"I2C_HIGH_SPPED_MODE\"\n$find.Replacement.Text = \"I2C_HIGH_SPEED_MODE\"\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll)\n"}
